# "error solve ifrs list": the per-share/aggregate figures in this sheet had been
# populated with the wrong magnitude (apparently whole-company totals instead of the
# correct per-unit figures for this ticker). Rewrite rows 2-6 (FY2014-FY2018 actuals)
# with the corrected figures, and blank out the forecast rows 7-9 (FY2019E-FY2021E),
# which turned out to have no reliable estimates and should only keep their
# period label (columns A-C) going forward.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4496
$ws.Range("E2").Value = 292
$ws.Range("F2").Value = 292
$ws.Range("G2").Value = 315
$ws.Range("H2").Value = 202
$ws.Range("I2").Value = 166
$ws.Range("J2").Value = 36
$ws.Range("K2").Value = 8586
$ws.Range("L2").Value = 1182
$ws.Range("M2").Value = 7404
$ws.Range("N2").Value = 6620
$ws.Range("O2").Value = 784
$ws.Range("P2").Value = 699
$ws.Range("Q2").Value = 619
$ws.Range("R2").Value = -88
$ws.Range("S2").Value = -341
$ws.Range("T2").Value = 19
$ws.Range("U2").Value = 600
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 4.48
$ws.Range("Y2").Value = 2.54
$ws.Range("Z2").Value = 2.35
$ws.Range("AA2").Value = 15.97
$ws.Range("AB2").Value = 837.96
$ws.Range("AC2").Value = 118
$ws.Range("AD2").Value = 33.82
$ws.Range("AE2").Value = 4732
$ws.Range("AF2").Value = 0.85
$ws.Range("AG2").Value = 25
$ws.Range("AH2").Value = 0.62
$ws.Range("AI2").Value = 21.11
$ws.Range("AJ2").Value = 139896190
$ws.Range("V2").ClearContents()

# Row 3
$ws.Range("D3").Value = 4735
$ws.Range("E3").Value = 233
$ws.Range("F3").Value = 233
$ws.Range("G3").Value = 401
$ws.Range("H3").Value = 252
$ws.Range("I3").Value = 225
$ws.Range("J3").Value = 27
$ws.Range("K3").Value = 8816
$ws.Range("L3").Value = 1224
$ws.Range("M3").Value = 7592
$ws.Range("N3").Value = 6795
$ws.Range("O3").Value = 797
$ws.Range("P3").Value = 699
$ws.Range("Q3").Value = 145
$ws.Range("R3").Value = -537
$ws.Range("S3").Value = -50
$ws.Range("T3").Value = 35
$ws.Range("U3").Value = 110
$ws.Range("V3").Value = 3
$ws.Range("W3").Value = 4.91
$ws.Range("X3").Value = 5.33
$ws.Range("Y3").Value = 3.36
$ws.Range("Z3").Value = 2.9
$ws.Range("AA3").Value = 16.12
$ws.Range("AB3").Value = 860.25
$ws.Range("AC3").Value = 161
$ws.Range("AD3").Value = 24.14
$ws.Range("AE3").Value = 4857
$ws.Range("AF3").Value = 0.8
$ws.Range("AG3").Value = 25
$ws.Range("AH3").Value = 0.64
$ws.Range("AI3").Value = 15.51
$ws.Range("AJ3").Value = 139896190

# Row 4
$ws.Range("D4").Value = 5045
$ws.Range("E4").Value = 214
$ws.Range("F4").Value = 214
$ws.Range("G4").Value = 320
$ws.Range("H4").Value = 133
$ws.Range("I4").Value = 126
$ws.Range("J4").Value = 7
$ws.Range("K4").Value = 9058
$ws.Range("L4").Value = 1389
$ws.Range("M4").Value = 7669
$ws.Range("N4").Value = 6878
$ws.Range("O4").Value = 791
$ws.Range("P4").Value = 699
$ws.Range("Q4").Value = 460
$ws.Range("R4").Value = -304
$ws.Range("S4").Value = 14
$ws.Range("T4").Value = 21
$ws.Range("U4").Value = 439
$ws.Range("V4").Value = 97
$ws.Range("W4").Value = 4.24
$ws.Range("X4").Value = 2.63
$ws.Range("Y4").Value = 1.85
$ws.Range("Z4").Value = 1.49
$ws.Range("AA4").Value = 18.11
$ws.Range("AB4").Value = 872.16
$ws.Range("AC4").Value = 90
$ws.Range("AD4").Value = 29.66
$ws.Range("AE4").Value = 4917
$ws.Range("AF4").Value = 0.54
$ws.Range("AG4").Value = 25
$ws.Range("AH4").Value = 0.93
$ws.Range("AI4").Value = 27.72
$ws.Range("AJ4").Value = 139896190

# Row 5
$ws.Range("D5").Value = 4219
$ws.Range("E5").Value = -43
$ws.Range("F5").Value = -43
$ws.Range("G5").Value = 139
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 48
$ws.Range("J5").Value = -43
$ws.Range("K5").Value = 8684
$ws.Range("L5").Value = 1104
$ws.Range("M5").Value = 7580
$ws.Range("N5").Value = 6932
$ws.Range("O5").Value = 648
$ws.Range("P5").Value = 699
$ws.Range("Q5").Value = 53
$ws.Range("R5").Value = -145
$ws.Range("S5").Value = -187
$ws.Range("T5").Value = 14
$ws.Range("U5").Value = 39
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = -1.02
$ws.Range("X5").Value = 0.12
$ws.Range("Y5").Value = 0.7
$ws.Range("Z5").Value = 0.06
$ws.Range("AA5").Value = 14.57
$ws.Range("AB5").Value = 882.84
$ws.Range("AC5").Value = 34
$ws.Range("AD5").Value = 88.84999999999999
$ws.Range("AE5").Value = 4955
$ws.Range("AF5").Value = 0.62
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 139896190
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("D6").Value = 3911
$ws.Range("E6").Value = 71
$ws.Range("F6").Value = 71
$ws.Range("G6").Value = 262
$ws.Range("H6").Value = 150
$ws.Range("I6").Value = 114
$ws.Range("K6").Value = 8802
$ws.Range("L6").Value = 1055
$ws.Range("M6").Value = 7747
$ws.Range("N6").Value = 7065
$ws.Range("P6").Value = 699
$ws.Range("Q6").Value = 351
$ws.Range("R6").Value = -67
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 29
$ws.Range("U6").Value = 322
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 1.82
$ws.Range("X6").Value = 3.83
$ws.Range("Y6").Value = 1.62
$ws.Range("Z6").Value = 1.71
$ws.Range("AA6").Value = 13.62
$ws.Range("AB6").Value = 902.38
$ws.Range("AC6").Value = 81
$ws.Range("AD6").Value = 25.13
$ws.Range("AE6").Value = 5050
$ws.Range("AF6").Value = 0.4
$ws.Range("AJ6").Value = 139896190
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()

# Rows 7-9: clear all data columns except A, B, C (per diff)
$ws.Range("D7:AJ9").ClearContents()
